$p = $ppt.ActivePresentation

# Slide 7: "Pokrok umelej inteligencie vo svete sachu"
# Merge the split runs (Garrym/Kasparovom/Deep/Blue/brute/force/enginov/AlphaZero/enginu)
# back into single runs per paragraph, and drop the stray trailing endParaRPr.
#
# NOTE: this host's TextRange.Text setter tries to reuse/diff against the
# previous run(s) to keep formatting; if the new text shares a literal
# prefix with the old first run it can leave a stray split run behind, and
# replacing only part of a multi-run paragraph can silently no-op. Priming
# the whole text frame with an unrelated placeholder value first avoids
# both problems, and operating on the whole shape (not per-paragraph)
# avoids it for the endParaRPr as well.
$s7 = $p.Slides.Item(7)
$shape7 = $s7.Shapes.Item(2)
$shape7.TextFrame.TextRange.Text = "1"
$shape7.TextFrame.TextRange.Text = "V roku 1997 umelá inteligencia prvý krát porazila človeka v šachu`r" + `
"Bol to zápas medzi vtedajším majstrom sveta Garrym Kasparovom a počítačom Deep Blue od spoločnosti IBM`r" + `
"Počítač využíval hrubú silu (brute force) tým že za ním nebola priveľká logika, snažil sa iba vypočítať všetky možné ťahy a prehľadával milióny pozícií`r" + `
"Od vtedy sa výkon superpočítačov posunul tak, že žiadny človek nemá šancu poraziť jeden z najlepších šachových enginov AlphaZero`r" + `
"V dnešnej dobe sami hráči študujú partie tohto enginu aby mohli vymyslieť nové stratégie"

# Slide 8: "Odhalovanie podvodov na sachovych turnajoch pomocou umelej inteligencie"
# Merge the "... daného " + "enginu" runs into a single run and drop the
# trailing endParaRPr.
$s8 = $p.Slides.Item(8)
$shape8 = $s8.Shapes.Item(2)
$shape8.TextFrame.TextRange.Text = "1"
$shape8.TextFrame.TextRange.Text = "Predpokladalo sa, že ľudia budú využívať umelú inteligenciu na to, aby sa naučili lepšie hrať šach, bohužiaľ je táto možnosť zneužívaná na podvádzanie`r" + `
"Umelá inteligencia sa používa na zistenie, či niektorí súťažiaci hrajú lepšie, ako by mali, vzhľadom na ich históriu hier`r" + `
"Napríklad stránka Chess.com používa svoju obrovskú databázu a porovnáva pravdepodobnosť hráča urobiť najlepší možný ťah podľa daného enginu"
